$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

# Data for new rows: rowNum, dateSerial, B, C, D
$data = @(
    ,@(302, 44376, 0, 2, 21.54243860404998)
    ,@(303, 44377, 0, 1, 10.77121930202499)
    ,@(304, 44378, 0, 1, 10.77121930202499)
    ,@(305, 44379, 0, 1, 10.77121930202499)
    ,@(306, 44380, 0, 1, 10.77121930202499)
    ,@(307, 44381, 0, 0, 0)
    ,@(308, 44382, 0, 0, 0)
    ,@(309, 44383, 0, 0, 0)
    ,@(310, 44384, 0, 0, 0)
    ,@(311, 44385, 0, 0, 0)
    ,@(312, 44386, 0, 0, 0)
    ,@(313, 44387, 0, 0, 0)
    ,@(314, 44388, 0, 0, 0)
    ,@(315, 44389, 0, 0, 0)
    ,@(316, 44390, 0, 0, 0)
    ,@(317, 44391, 1, 1, 10.77121930202499)
    ,@(318, 44392, 1, 2, 21.54243860404998)
    ,@(319, 44393, 0, 2, 21.54243860404998)
    ,@(320, 44394, 1, 3, 32.31365790607497)
    ,@(321, 44395, 0, 3, 32.31365790607497)
    ,@(322, 44396, 0, 3, 32.31365790607497)
    ,@(323, 44397, 0, 3, 32.31365790607497)
    ,@(324, 44398, 0, 2, 21.54243860404998)
    ,@(325, 44399, 0, 1, 10.77121930202499)
    ,@(326, 44400, 2, 3, 32.31365790607497)
    ,@(327, 44401, 0, 2, 21.54243860404998)
    ,@(328, 44402, 2, 4, 43.08487720809995)
)

$lastExistingRow = 301

foreach ($entry in $data) {
    $rowNum = $entry[0]
    $dateVal = $entry[1]
    $bVal = $entry[2]
    $cVal = $entry[3]
    $dVal = $entry[4]

    $srcRange = $ws.Range("A" + $lastExistingRow + ":D" + $lastExistingRow)
    $dstRange = $ws.Range("A" + $rowNum + ":D" + $rowNum)
    $srcRange.Copy($dstRange)

    $ws.Cells.Item($rowNum, 1).Value = $dateVal
    $ws.Cells.Item($rowNum, 2).Value = $bVal
    $ws.Cells.Item($rowNum, 3).Value = $cVal
    $ws.Cells.Item($rowNum, 4).Value = $dVal
}

Write-Host "Done adding rows 302-328"
